$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12
$ws.Range("G12").Value = 1240524717.0500007
$ws.Range("I12").Value = 4188377156

# Row 13
$ws.Range("G13").Value = 319819483.18000001
$ws.Range("I13").Value = 1012006300

# Row 14
$ws.Range("G14").Value = 34063116.800000042
$ws.Range("I14").Value = 44319159.289999999

# Row 16
$ws.Range("I16").Value = 162861893.59999999

# Row 18 - G18 becomes a formula
$ws.Range("G18").Formula = "=SUM(G12:G17)"

# Row 19
$ws.Range("G19").Value = -379300000.00000012
$ws.Range("I19").Value = 1160500000

# Row 21 - G21 becomes a formula
$ws.Range("G21").Formula = "=SUM(G18:G20)"

# Row 22
$ws.Range("I22").Value = 85592745

# Row 26
$ws.Range("G26").Value = 1029174575.116062
$ws.Range("I26").Value = 1010658959

$wb.Application.CalculateFullRebuild()
